$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 41
$ws.Range("D4").Value = 47
$ws.Range("D5").Value = 46
$ws.Range("J5").Value = 21
